$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from serial date
# 45183 (2023-09-14) to 45184 (2023-09-15).
foreach ($r in 2..6) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
